$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the "time_laptop" tab to "energy_intensity_network" and make it
#    the active sheet (activeTab moves from 0 -> 1, tabSelected moves off
#    "params" and onto this sheet).
# ---------------------------------------------------------------------------
$wsEnergy = $wb.Worksheets.Item("time_laptop")
$wsEnergy.Name = "energy_intensity_network"

$wsParams = $wb.Worksheets.Item("params")

# ---------------------------------------------------------------------------
# 2. "params" sheet: rows 2 and 6 swap their entire contents.
#    Row 2 was carbon_intensity / exp ...   -> becomes power_latop / interp ...
#    Row 6 was power_latop / interp ...     -> becomes carbon_intensity / exp ...
#    (the ref-date column I is identical in both rows, so it is left alone)
# ---------------------------------------------------------------------------

# --- new row 2 (previously row 6: power_latop / interp) ---
$wsParams.Range("A2").Value = "power_latop"
$wsParams.Range("C2").Value = "interp"
$wsParams.Range("D2").Value = "linear"
$wsParams.Range("E2").Value = "{""2020-01-01"":10, ""2031-06-01"":9.5}"
$wsParams.Range("F2").Value = 0
$wsParams.Range("G2").Value = 4
$wsParams.Range("H2").Value = 0.05
$wsParams.Range("J2").Value = "W"
$wsParams.Range("P2").Value = "what does it mean? How do collect this info?"
$wsParams.Range("Q2").Value = "x"
$wsParams.Range("R2").Value = "power draw of laptop"
$wsParams.Range("S2").Value = 0

# --- new row 6 (previously row 2: carbon_intensity / exp) ---
$wsParams.Range("A6").Value = "carbon_intensity"
$wsParams.Range("C6").Value = "exp"
$wsParams.Range("D6").Clear()
$wsParams.Range("E6").Value = 0.5
$wsParams.Range("F6").Value = -0.1
$wsParams.Range("G6").Value = 0.1
$wsParams.Range("H6").Value = 0.05
$wsParams.Range("J6").Value = "kg/kWh"
$wsParams.Range("P6").Clear()
$wsParams.Range("Q6").Clear()
$wsParams.Range("R6").Clear()
$wsParams.Range("S6").Value = 5

# "params" selection becomes a single cell, A4, and loses tabSelected (handled
# automatically once another sheet is activated below).
$wsParams.Range("A4").Select()

# ---------------------------------------------------------------------------
# 3. "energy_intensity_network" sheet: update the UK / DE rows.
# ---------------------------------------------------------------------------
$wsEnergy.Range("C2").Value = 20
$wsEnergy.Range("D2").Value = 0.25
$wsEnergy.Range("E2").Value = 0.1
$wsEnergy.Range("F2").Value = 0.1
$wsEnergy.Range("G2").Value = 6

# G2 previously carried a custom "integer" number format (style index 7); the
# edit drops that custom formatting, leaving G2 with the plain default style
# - same as its neighbour G3. Copy G3's (default) formatting onto G2.
$wsEnergy.Range("G3").Copy()
$wsEnergy.Range("G2").PasteSpecial(-4122)

$wsEnergy.Range("C3").Value = 20
$wsEnergy.Range("D3").Value = 0.25
$wsEnergy.Range("E3").Value = 0.1
$wsEnergy.Range("F3").Value = 0.1
$wsEnergy.Range("G3").Value = 7

# Activate this sheet last so it becomes the workbook's active tab, with the
# whole table selected (A1:G3, active cell A1).
$wsEnergy.Activate()
$wsEnergy.Range("A1:G3").Select()
